$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1015.7143
$ws.Range("I12").Value = 1277.75
$ws.Range("J12").Value = 666.3333
$ws.Range("K12").Value = 1277.75
$ws.Range("L12").Value = 666.3333
$ws.Range("M12").Value = -1107.75
$ws.Range("N12").Value = -1006.3333
$ws.Range("H38").Value = 5059.3335
$ws.Range("I38").Value = 264.5
$ws.Range("K38").Value = 793.5
$ws.Range("M38").Value = -421.5
$ws.Range("H76").Value = 5142.4707
$ws.Range("I76").Value = 4428
$ws.Range("J76").Value = 5777.5557
$ws.Range("K76").Value = 4428
$ws.Range("L76").Value = 5777.5557
$ws.Range("M76").Value = -4113
$ws.Range("N76").Value = -6407.5557
$ws.Range("H79").Value = 5142.4707
$ws.Range("I79").Value = 4428
$ws.Range("J79").Value = 5777.5557
$ws.Range("K79").Value = 4428
$ws.Range("L79").Value = 5777.5557
$ws.Range("M79").Value = -3336
$ws.Range("N79").Value = -7961.5557
$ws.Range("H106").Value = 2434.7222
$ws.Range("I106").Value = 2869.25
$ws.Range("K106").Value = 2869.25
$ws.Range("M106").Value = -2238.25
$ws.Range("H116").Value = 20241112
$ws.Range("I116").Value = 26987134
$ws.Range("J116").Value = 3040.5715
$ws.Range("K116").Value = 26987134
$ws.Range("L116").Value = 3040.5715
$ws.Range("M116").Value = -26983692
$ws.Range("N116").Value = -9924.5715
$ws.Range("H132").Value = 6162.86
$ws.Range("J132").Value = 27770.445
$ws.Range("L132").Value = 83311.33499999999
$ws.Range("N132").Value = -88371.33499999999
$ws.Range("H138").Value = 6240.245
$ws.Range("J138").Value = 6912.298
$ws.Range("L138").Value = 20736.894
$ws.Range("N138").Value = -31016.894

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 16400
$ws.Range("I39").Value = 16400
$ws.Range("K39").Value = 16400
$ws.Range("M39").Value = -15880
$ws.Range("H45").Value = 55835.332
$ws.Range("I45").Value = 64196.777
$ws.Range("J45").Value = 5666.6665
$ws.Range("K45").Value = 64196.777
$ws.Range("L45").Value = 5666.6665
$ws.Range("M45").Value = -63819.777
$ws.Range("N45").Value = -6420.6665
$ws.Range("H61").Value = 4288.026
$ws.Range("I61").Value = 4128.6943
$ws.Range("K61").Value = 4128.6943
$ws.Range("M61").Value = -3916.6943
$ws.Range("H63").Value = 1997
$ws.Range("I63").Value = 1997
$ws.Range("K63").Value = 1997
$ws.Range("M63").Value = -1311
$ws.Range("H66").Value = 1997
$ws.Range("I66").Value = 1997
$ws.Range("K66").Value = 9985
$ws.Range("M66").Value = -6553
$ws.Range("H74").Value = 976.2353000000001
$ws.Range("I74").Value = 916
$ws.Range("J74").Value = 1940
$ws.Range("K74").Value = 916
$ws.Range("L74").Value = 1940
$ws.Range("M74").Value = -42
$ws.Range("N74").Value = -3688
$ws.Range("H77").Value = 976.2353000000001
$ws.Range("I77").Value = 916
$ws.Range("J77").Value = 1940
$ws.Range("K77").Value = 4580
$ws.Range("L77").Value = 9700
$ws.Range("M77").Value = -212
$ws.Range("N77").Value = -18436
$ws.Range("H102").Value = 1204.1364
$ws.Range("I102").Value = 1192.4117
$ws.Range("J102").Value = 1244
$ws.Range("K102").Value = 1192.4117
$ws.Range("L102").Value = 1244
$ws.Range("M102").Value = 429.5882999999999
$ws.Range("N102").Value = -4488
$ws.Range("H110").Value = 1528.4445
$ws.Range("I110").Value = 1322.7142
$ws.Range("K110").Value = 1322.7142
$ws.Range("M110").Value = 722.2858000000001
$ws.Range("H136").Value = 4288.026
$ws.Range("I136").Value = 4128.6943
$ws.Range("K136").Value = 12386.0829
$ws.Range("M136").Value = -9836.082900000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 15000
$ws.Range("I33").Value = 15000
$ws.Range("K33").Value = 15000
$ws.Range("M33").Value = -14664
$ws.Range("H99").Value = 2754.8235
$ws.Range("I99").Value = 2922.2666
$ws.Range("K99").Value = 2922.2666
$ws.Range("M99").Value = -1424.2666
$ws.Range("H139").Value = 110008.8
$ws.Range("I139").Value = 106779
$ws.Range("J139").Value = 110367.664
$ws.Range("K139").Value = 106779
$ws.Range("L139").Value = 110367.664
$ws.Range("N139").Value = -120647.664
$ws.Range("M139").Value = -101639

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 152.86667
$ws.Range("J7").Value = 144.14285
$ws.Range("L7").Value = 144.14285
$ws.Range("N7").Value = -370.14285
$ws.Range("H31").Value = 21742554
$ws.Range("J31").Value = 7355.8
$ws.Range("L31").Value = 7355.8
$ws.Range("N31").Value = -7945.8
$ws.Range("H34").Value = 21742554
$ws.Range("J34").Value = 7355.8
$ws.Range("L34").Value = 7355.8
$ws.Range("N34").Value = -7759.8
$ws.Range("H58").Value = 2964.9583
$ws.Range("I58").Value = 2564.5334
$ws.Range("J58").Value = 3632.3333
$ws.Range("K58").Value = 2564.5334
$ws.Range("L58").Value = 3632.3333
$ws.Range("M58").Value = -2361.5334
$ws.Range("N58").Value = -4038.3333
$ws.Range("H86").Value = 5391.0435
$ws.Range("I86").Value = 5876.3
$ws.Range("K86").Value = 5876.3
$ws.Range("M86").Value = -4753.3
$ws.Range("H89").Value = 5391.0435
$ws.Range("I89").Value = 5876.3
$ws.Range("K89").Value = 29381.5
$ws.Range("M89").Value = -23765.5
$ws.Range("H132").Value = 34189664
$ws.Range("I132").Value = 39217156
$ws.Range("K132").Value = 117651468
$ws.Range("M132").Value = -117648938
$ws.Range("H134").Value = 2592.7827
$ws.Range("I134").Value = 1655.2941
$ws.Range("K134").Value = 4965.8823
$ws.Range("M134").Value = -2430.8823
$ws.Range("H136").Value = 2964.9583
$ws.Range("I136").Value = 2564.5334
$ws.Range("J136").Value = 3632.3333
$ws.Range("K136").Value = 7693.600199999999
$ws.Range("L136").Value = 10896.9999
$ws.Range("M136").Value = -5143.600199999999
$ws.Range("N136").Value = -15996.9999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 421.25
$ws.Range("I51").Value = 421.25
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 1263.75
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -803.75
$ws.Range("N51").ClearContents()
$ws.Range("H81").Value = 3952.4614
$ws.Range("J81").Value = 4787.75
$ws.Range("L81").Value = 14363.25
$ws.Range("N81").Value = -16609.25
$ws.Range("H84").Value = 3952.4614
$ws.Range("J84").Value = 4787.75
$ws.Range("L84").Value = 43089.75
$ws.Range("N84").Value = -54321.75
$ws.Range("H86").Value = 1199.2
$ws.Range("I86").Value = 1666
$ws.Range("J86").Value = 499
$ws.Range("K86").Value = 4998
$ws.Range("L86").Value = 1497
$ws.Range("M86").Value = -3812
$ws.Range("N86").Value = -3869
$ws.Range("H89").Value = 1199.2
$ws.Range("I89").Value = 1666
$ws.Range("J89").Value = 499
$ws.Range("K89").Value = 14994
$ws.Range("L89").Value = 4491
$ws.Range("M89").Value = -9066
$ws.Range("N89").Value = -16347
$ws.Range("H111").Value = 6615.875
$ws.Range("I111").Value = 4919.857
$ws.Range("K111").Value = 14759.571
$ws.Range("M111").Value = -11692.571
$ws.Range("H113").Value = 678.9091
$ws.Range("J113").Value = 728.2857
$ws.Range("L113").Value = 2184.8571
$ws.Range("N113").Value = -6524.8571
$ws.Range("H121").Value = 1143.3636
$ws.Range("I121").Value = 299
$ws.Range("K121").Value = 897
$ws.Range("M121").Value = 413

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 120199.9
$ws.Range("I80").Value = 207000
$ws.Range("J80").Value = 33399.8
$ws.Range("K80").Value = 207000
$ws.Range("L80").Value = 33399.8
$ws.Range("M80").Value = -206002
$ws.Range("N80").Value = -35395.8
$ws.Range("H83").Value = 120199.9
$ws.Range("I83").Value = 207000
$ws.Range("J83").Value = 33399.8
$ws.Range("K83").Value = 1035000
$ws.Range("L83").Value = 166999
$ws.Range("M83").Value = -1030008
$ws.Range("N83").Value = -176983
$ws.Range("H102").Value = 558230.9
$ws.Range("I102").Value = 662253.1
$ws.Range("K102").Value = 662253.1
$ws.Range("M102").Value = -660631.1
$ws.Range("H126").Value = 1652.9333
$ws.Range("I126").Value = 1652.9333
$ws.Range("K126").Value = 4958.7999
$ws.Range("M126").Value = -2488.7999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 12402.833
$ws.Range("I122").Value = 5055.75
$ws.Range("K122").Value = 15167.25
$ws.Range("M122").Value = -12717.25
$ws.Range("H132").Value = 2659.34
$ws.Range("I132").Value = 2598.7715
$ws.Range("K132").Value = 7796.314499999999
$ws.Range("M132").Value = -5266.314499999999
$ws.Range("H136").Value = 4284.907
$ws.Range("I136").Value = 2794.6128
$ws.Range("K136").Value = 8383.838400000001
$ws.Range("M136").Value = -5833.838400000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1014.2727
$ws.Range("I107").Value = 665.2857
$ws.Range("K107").Value = 1995.8571
$ws.Range("M107").Value = -75.85710000000017
$ws.Range("H113").Value = 1399.8182
$ws.Range("J113").Value = 1789.8
$ws.Range("L113").Value = 5369.4
$ws.Range("N113").Value = -9709.4
$ws.Range("H136").Value = 3643.1516
$ws.Range("I136").Value = 2579.6
$ws.Range("K136").Value = 7738.799999999999
$ws.Range("M136").Value = -5188.799999999999
